$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source publishes a new weekly price observation (Primera/Segunda
# quality rows) for "Betarraga" on top of the historical series. Insert two
# fresh rows at the top of the data block (352:353), which pushes the whole
# existing history down by two rows (old 352-375 -> new 354-377).
$ws.Range("A352:R353").EntireRow.Insert()

# The newly inserted rows are blank; populate the columns that stay
# constant across observations by copying them from the row immediately
# below (which now holds what used to be row 352 before the insert).
$constantCols = @("A", "B", "C", "E", "F", "G", "H", "I", "N", "O", "Q", "R")
foreach ($col in $constantCols) {
    $ws.Range($col + "352").Value = $ws.Range($col + "354").Value2
    $ws.Range($col + "353").Value = $ws.Range($col + "355").Value2
}

# Quality differs between the two new rows, same as every other pair.
$ws.Range("I352").Value = "Primera"
$ws.Range("I353").Value = "Segunda"
$ws.Range("N352").Value = "$/paquete 4 unidades"
$ws.Range("N353").Value = "$/paquete 5 unidades"
$ws.Range("Q352").Value = 4
$ws.Range("Q353").Value = 5

# Now write the new week's figures for the two quality rows.
$ws.Range("D352").Value = 44826
$ws.Range("J352").Value = 800
$ws.Range("K352").Value = 600
$ws.Range("L352").Value = 700
$ws.Range("M352").Value = 650
$ws.Range("P352").Value = 162

$ws.Range("D353").Value = 44826
$ws.Range("J353").Value = 800
$ws.Range("K353").Value = 600
$ws.Range("L353").Value = 700
$ws.Range("M353").Value = 650
$ws.Range("P353").Value = 130
